$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8th column). Excel copies the
# left-neighbour column's (G) formatting onto the freshly inserted column,
# and shifts every later column (old H..M) one slot to the right (I..N).
$ws.Columns.Item(8).EntireColumn.Insert()

# Match the new column's width to its left neighbour (column G) - same as
# what the inserted column already visually inherited.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Give the new column its header label.
$ws.Range("H6").Value = "Phát sinh tăng"

# The "Tổng:" totals row gets a dedicated bordered style for the new column
# (top/bottom border only, no left/right) instead of the plain style it
# inherited from its neighbours.
$ws.Range("H31").Borders.Item(7).LineStyle = 0    # xlEdgeLeft -> none
$ws.Range("H31").Borders.Item(10).LineStyle = 0   # xlEdgeRight -> none
$ws.Range("H31").Borders.Item(8).LineStyle = 1    # xlEdgeTop -> continuous
$ws.Range("H31").Borders.Item(9).LineStyle = 1    # xlEdgeBottom -> continuous

# Restore the active-cell selection to match the edited sheet's cursor spot.
[void]$ws.Range("B12").Select()
